$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new row (22) for the "Bili video link parsing" command, below the
#     existing last row (21, "签到"). We mirror row 21's layout: col A holds
#     the command name, col H holds the version tag, B..G stay blank. ---

$ws.Range("A22").Value = "Bili视频链接解析"
$ws.Range("H22").Value = "0.0.1 Pre"

$ws.Rows("22:22").RowHeight = $ws.Rows("21:21").RowHeight

# Copy the visual formatting of row 21 onto row 22, cell by cell, so the new
# row matches the look of the rest of the table (font, fill, alignment,
# borders on all four sides).
$srcRow = $ws.Range("A21:H21")
$dstRow = $ws.Range("A22:H22")

for ($col = 1; $col -le 8; $col++) {
    $srcCell = $srcRow.Cells.Item(1, $col)
    $dstCell = $dstRow.Cells.Item(1, $col)

    $dstCell.Font.Name = $srcCell.Font.Name
    $dstCell.Font.Size = $srcCell.Font.Size
    $dstCell.Font.Bold = $srcCell.Font.Bold
    $dstCell.Font.Color = $srcCell.Font.Color

    $dstCell.Interior.Pattern = $srcCell.Interior.Pattern
    $dstCell.Interior.Color = $srcCell.Interior.Color

    $dstCell.HorizontalAlignment = $srcCell.HorizontalAlignment
    $dstCell.VerticalAlignment = $srcCell.VerticalAlignment

    foreach ($edge in 7, 8, 9, 10) {
        $sb = $srcCell.Borders.Item($edge)
        $db = $dstCell.Borders.Item($edge)
        $db.LineStyle = $sb.LineStyle
        if ($sb.LineStyle -ne -4142) {
            $db.Weight = $sb.Weight
            $db.Color = $sb.Color
        }
    }
}

# Sheet used to end at row 21 / the selection handle covered A1:H22 — now
# that row 22 is real data, the dimension grows and the leftover selection
# marker is replaced with a normal single-cell selection further down.
$ws.Range("F25").Select()

# Command column was a bit tight for the longer new label - widen it.
$ws.Columns("A").ColumnWidth = 18.5
